$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing row 16, column A value (trim floating point noise)
$ws.Range("A16").Value = 45876.58349

# Append new row 17 with the latest sensor reading
$ws.Range("A17").Value = 45876.62519435577
$ws.Range("B17").Value = 2025
$ws.Range("C17").Value = 28
$ws.Range("D17").Value = 19.06
$ws.Range("E17").Value = 78.55
$ws.Range("F17").Value = 450
$ws.Range("G17").Value = 16.75
$ws.Range("H17").Value = "ESE"
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = "15:00:16"

# Match the date/time number formatting used by the rest of column A
$ws.Range("A17").NumberFormat = $ws.Range("A16").NumberFormat
